# Apply a permutation of full-row contents (columns A:AY) for the rows
# affected by the source commit. Each row in $mapping takes on the
# content that currently belongs to the row given as its value; the
# permutation is closed over the same set of rows, so we must snapshot
# every row's values BEFORE writing any of them back.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destinationRow -> sourceRow (destination receives source's current content)
$mapping = @{
    115 = 116
    116 = 117
    117 = 118
    118 = 119
    119 = 120
    120 = 121
    121 = 122
    122 = 123
    123 = 115
    124 = 130
    125 = 124
    126 = 125
    127 = 126
    128 = 127
    129 = 128
    130 = 129
    136 = 140
    137 = 136
    138 = 137
    139 = 138
    140 = 139
    159 = 162
    160 = 163
    161 = 164
    162 = 159
    163 = 160
    164 = 161
    174 = 185
    175 = 174
    176 = 175
    177 = 176
    178 = 177
    179 = 178
    180 = 179
    181 = 180
    182 = 181
    183 = 182
    184 = 183
    185 = 184
}

$firstCol = "A"
$lastCol = "AY"

# Snapshot every involved row's current values first.
$snapshot = @{}
foreach ($r in $mapping.Keys) {
    $addr = "$firstCol$r" + ":" + "$lastCol$r"
    $snapshot[$r] = $ws.Range($addr).Value()
}

# Now write each destination row the snapshot taken from its source row.
foreach ($r in $mapping.Keys) {
    $src = $mapping[$r]
    $addr = "$firstCol$r" + ":" + "$lastCol$r"
    $ws.Range($addr).Value = $snapshot[$src]
}
